$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.437.47'
$ws.Range("E2").Value = '  -0.47%  '

$ws.Range("D3").Value = '3.496.19'
$ws.Range("E3").Value = '  -0.60%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.76'
$ws.Range("E5").Value = '  -2.87%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '193.60'
$ws.Range("E6").Value = '  +0.60%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.618'
$ws.Range("E7").Value = '  -1.92%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("E9").Value = '  -6.29%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.642'
$ws.Range("E10").Value = '  -3.52%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.84'
$ws.Range("E11").Value = '  -1.23%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000296'
$ws.Range("E12").Value = '  -4.10%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.40'
$ws.Range("E13").Value = '  -2.07%  '

$ws.Range("D14").Value = '4.057.06'
$ws.Range("E14").Value = '  -0.67%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '595.03'
$ws.Range("E15").Value = '  -4.28%  '

$ws.Range("D16").Value = '69.582.33'
$ws.Range("E16").Value = '  -0.41%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.86'
$ws.Range("E17").Value = '  -0.94%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.57'
$ws.Range("E18").Value = '  -0.67%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.122'
$ws.Range("E19").Value = '  +2.01%  '

$ws.Range("D20").Value = '3.498.00'
$ws.Range("E20").Value = '  -0.29%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.976'
$ws.Range("E21").Value = '  -1.56%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.72'
$ws.Range("E22").Value = '  +3.05%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.24'
$ws.Range("E23").Value = '  +3.62%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '101.99'
$ws.Range("E24").Value = '  -6.92%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.57'
$ws.Range("E25").Value = '  -3.09%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.08'
$ws.Range("E26").Value = '  -1.14%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.74'
$ws.Range("E27").Value = '  -2.57%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.41'
$ws.Range("E28").Value = '  -3.04%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '32.84'
$ws.Range("E29").Value = '  -4.04%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.22'
$ws.Range("E30").Value = '  +7.45%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.92'
$ws.Range("E31").Value = '  -1.04%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.23'
$ws.Range("E32").Value = '  -2.41%  '

$ws.Range("E33").Value = '  -2.93%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.09'
$ws.Range("E34").Value = '  -0.54%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.19'
$ws.Range("E35").Value = '  +2.26%  '

$ws.Range("D36").Value = '3.736.95'
$ws.Range("E36").Value = '  +2.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'

$ws.Range("D38").Value = '0.0₃0803'
$ws.Range("E38").Value = '  +2.86%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.63'
$ws.Range("E39").Value = '  -0.57%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.386'
$ws.Range("E40").Value = '  -2.17%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '494.18'
$ws.Range("E41").Value = '  -4.48%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '35.74'
$ws.Range("E42").Value = '  -2.68%  '

$ws.Range("E43").Value = '  -4.62%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0446'
$ws.Range("E44").Value = '  -5.34%  '

$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.138'
$ws.Range("E45").Value = '  -3.26%  '

$ws.Range("B46").Value = 'ThetaToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.79'
$ws.Range("E46").Value = '  -4.81%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.22'
$ws.Range("E47").Value = '  -3.41%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.01'
$ws.Range("E48").Value = '  +0.18%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.36'
$ws.Range("E49").Value = '  -4.83%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.35'
$ws.Range("E50").Value = '  +0.47%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000241'
$ws.Range("E51").Value = '  +0.11%  '
